$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First pass: populate cells in the order needed to build the shared-string
# table in the same sequence as the target workbook (Excel appends new unique
# strings to the shared string table in first-use order). ---

$ws.Range("C1").Value = "Field of Purchase"
$ws.Range("C2").Value = "Electronics"
$ws.Range("C3").Value = "Cloths"
$ws.Range("C4").Value = "sports"
$ws.Range("C5").Value = "Delicacies"
$ws.Range("C6").Value = "OTT"

$ws.Range("D1").Value = "Country"
$ws.Range("D2").Value = "India"
$ws.Range("D3").Value = "china"
$ws.Range("D4").Value = "Canada"
$ws.Range("D5").Value = "US"
$ws.Range("D6").Value = "Mexico"

# --- Second pass: correct the few cells whose final value differs from the
# initial entry above, re-using already-registered shared strings. ---

$ws.Range("C3").Value = "Delicacies"
$ws.Range("C5").Value = "OTT"
$ws.Range("C6").Value = "Cloths"

# Autofit columns A, B, C to match bestFit widths seen in target
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Update selection to match target (E4)
$ws.Range("E4").Select() | Out-Null
